$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 138.9
$ws.Range("I11").Value = 138.9
$ws.Range("K11").Value = 138.9
$ws.Range("M11").Value = 1.099999999999994

$ws.Range("H21").Value = 14978.875
$ws.Range("I21").Value = 4017
$ws.Range("J21").Value = 18632.834
$ws.Range("K21").Value = 4017
$ws.Range("L21").Value = 18632.834
$ws.Range("M21").Value = -3549
$ws.Range("N21").Value = -19568.834

$ws.Range("H23").Value = 14978.875
$ws.Range("I23").Value = 4017
$ws.Range("J23").Value = 18632.834
$ws.Range("K23").Value = 4017
$ws.Range("L23").Value = 18632.834
$ws.Range("M23").Value = -3783
$ws.Range("N23").Value = -19100.834

$ws.Range("H29").Value = 180.75
$ws.Range("I29").Value = 180.75
$ws.Range("K29").Value = 542.25
$ws.Range("M29").Value = -261.25

$ws.Range("H33").Value = 336.7619
$ws.Range("I33").Value = 348.6
$ws.Range("K33").Value = 348.6
$ws.Range("M33").Value = -119.6

$ws.Range("H38").Value = 352.25
$ws.Range("I38").Value = 302.25
$ws.Range("J38").Value = 502.25
$ws.Range("K38").Value = 906.75
$ws.Range("L38").Value = 1506.75
$ws.Range("M38").Value = -534.75
$ws.Range("N38").Value = -2250.75

$ws.Range("H58").Value = 1499.9286
$ws.Range("I58").Value = 1119.9
$ws.Range("J58").Value = 2450
$ws.Range("K58").Value = 3359.7
$ws.Range("L58").Value = 7350
$ws.Range("M58").Value = -3209.7
$ws.Range("N58").Value = -7650

$ws.Range("H96").Value = 1054.6364
$ws.Range("I96").Value = 755.6667
$ws.Range("K96").Value = 2267.0001
$ws.Range("M96").Value = -894.0001000000002

$ws.Range("H127").Value = 869.55
$ws.Range("I127").Value = 433.9091
$ws.Range("J127").Value = 1402
$ws.Range("K127").Value = 1301.7273
$ws.Range("L127").Value = 4206
$ws.Range("M127").Value = 3658.2727
$ws.Range("N127").Value = -14126

$ws.Range("H137").Value = 2424.1406
$ws.Range("I137").Value = 2110.5813
$ws.Range("K137").Value = 6331.743899999999
$ws.Range("M137").Value = -3781.743899999999

$ws.Range("H138").Value = 2563.9849
$ws.Range("J138").Value = 3565.8125
$ws.Range("L138").Value = 10697.4375
$ws.Range("N138").Value = -20977.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1681.4
$ws.Range("I2").Value = 1743
$ws.Range("J2").Value = 1537.6666
$ws.Range("K2").Value = 1743
$ws.Range("L2").Value = 1537.6666
$ws.Range("M2").Value = -1630
$ws.Range("N2").Value = -1763.6666

$ws.Range("H74").Value = 2674.9473
$ws.Range("I74").Value = 2070
$ws.Range("J74").Value = 3985.6667
$ws.Range("K74").Value = 2070
$ws.Range("L74").Value = 3985.6667
$ws.Range("M74").Value = -1196
$ws.Range("N74").Value = -5733.6667

$ws.Range("H77").Value = 2674.9473
$ws.Range("I77").Value = 2070
$ws.Range("J77").Value = 3985.6667
$ws.Range("K77").Value = 10350
$ws.Range("L77").Value = 19928.3335
$ws.Range("M77").Value = -5982
$ws.Range("N77").Value = -28664.3335

$ws.Range("H110").Value = 1634.9231
$ws.Range("I110").Value = 1618.5454
$ws.Range("K110").Value = 1618.5454
$ws.Range("M110").Value = 426.4546

$ws.Range("H116").Value = 1681.4
$ws.Range("I116").Value = 1743
$ws.Range("J116").Value = 1537.6666
$ws.Range("K116").Value = 1743
$ws.Range("L116").Value = 1537.6666
$ws.Range("M116").Value = 551
$ws.Range("N116").Value = -6125.6666

$ws.Range("H132").Value = 3415.9534
$ws.Range("I132").Value = 2293.121
$ws.Range("J132").Value = 7121.3
$ws.Range("K132").Value = 6879.363
$ws.Range("L132").Value = 21363.9
$ws.Range("M132").Value = -4349.363
$ws.Range("N132").Value = -26423.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1681.4
$ws.Range("I3").Value = 1743
$ws.Range("J3").Value = 1537.6666
$ws.Range("K3").Value = 1743
$ws.Range("L3").Value = 1537.6666
$ws.Range("M3").Value = -1629
$ws.Range("N3").Value = -1765.6666

$ws.Range("H99").Value = 2550

$ws.Range("H105").Value = 6580724
$ws.Range("I105").Value = 12501410
$ws.Range("J105").Value = 2184.3333
$ws.Range("K105").Value = 12501410
$ws.Range("L105").Value = 2184.3333
$ws.Range("M105").Value = -12499663
$ws.Range("N105").Value = -5678.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6124.4443
$ws.Range("I31").Value = 1177.6552
$ws.Range("J31").Value = 11862.72
$ws.Range("K31").Value = 1177.6552
$ws.Range("L31").Value = 11862.72
$ws.Range("M31").Value = -882.6551999999999
$ws.Range("N31").Value = -12452.72

$ws.Range("H34").Value = 6124.4443
$ws.Range("I34").Value = 1177.6552
$ws.Range("J34").Value = 11862.72
$ws.Range("K34").Value = 1177.6552
$ws.Range("L34").Value = 11862.72
$ws.Range("M34").Value = -975.6551999999999
$ws.Range("N34").Value = -12266.72

$ws.Range("H58").Value = 1708.6522
$ws.Range("I58").Value = 1356.4
$ws.Range("J58").Value = 1979.6154
$ws.Range("K58").Value = 1356.4
$ws.Range("L58").Value = 1979.6154
$ws.Range("M58").Value = -1153.4
$ws.Range("N58").Value = -2385.6154

$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

$ws.Range("H132").Value = 5210048
$ws.Range("I132").Value = 1377.96
$ws.Range("J132").Value = 23812440
$ws.Range("K132").Value = 4133.88
$ws.Range("L132").Value = 71437320
$ws.Range("M132").Value = -1603.88
$ws.Range("N132").Value = -71442380

$ws.Range("H136").Value = 1708.6522
$ws.Range("I136").Value = 1356.4
$ws.Range("J136").Value = 1979.6154
$ws.Range("K136").Value = 4069.2
$ws.Range("L136").Value = 5938.8462
$ws.Range("M136").Value = -1519.2
$ws.Range("N136").Value = -11038.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 925.13336
$ws.Range("I5").Value = 587
$ws.Range("K5").Value = 1761
$ws.Range("M5").Value = -1649

$ws.Range("H62").Value = 8434.429
$ws.Range("J62").Value = 8434.429
$ws.Range("L62").Value = 25303.287
$ws.Range("N62").Value = -26675.287

$ws.Range("H65").Value = 8434.429
$ws.Range("J65").Value = 8434.429
$ws.Range("L65").Value = 75909.861
$ws.Range("N65").Value = -82773.861

$ws.Range("H68").Value = 842
$ws.Range("I68").Value = 1151
$ws.Range("J68").Value = 687.5
$ws.Range("K68").Value = 3453
$ws.Range("L68").Value = 2062.5
$ws.Range("M68").Value = -2642
$ws.Range("N68").Value = -3684.5

$ws.Range("H69").Value = 2575
$ws.Range("I69").Value = 494
$ws.Range("J69").Value = 2991.2
$ws.Range("K69").Value = 1482
$ws.Range("L69").Value = 8973.599999999999
$ws.Range("M69").Value = -671
$ws.Range("N69").Value = -10595.6

$ws.Range("H71").Value = 842
$ws.Range("I71").Value = 1151
$ws.Range("J71").Value = 687.5
$ws.Range("K71").Value = 10359
$ws.Range("L71").Value = 6187.5
$ws.Range("M71").Value = -6303
$ws.Range("N71").Value = -14299.5

$ws.Range("H72").Value = 2575
$ws.Range("I72").Value = 494
$ws.Range("J72").Value = 2991.2
$ws.Range("K72").Value = 4446
$ws.Range("L72").Value = 26920.8
$ws.Range("M72").Value = -390
$ws.Range("N72").Value = -35032.8

$ws.Range("H74").Value = 3900
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 11700
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -13822

$ws.Range("H75").Value = 1500
$ws.Range("I75").Value = 1500
$ws.Range("K75").Value = 4500
$ws.Range("M75").Value = -3502

$ws.Range("H77").Value = 3900
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 35100
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -45708

$ws.Range("H78").Value = 1500
$ws.Range("I78").Value = 1500
$ws.Range("K78").Value = 13500
$ws.Range("M78").Value = -8508

$ws.Range("H107").Value = 71429000
$ws.Range("J107").Value = 166667230
$ws.Range("L107").Value = 500001690
$ws.Range("N107").Value = -500005530

$ws.Range("H122").Value = 4832.24
$ws.Range("J122").Value = 7718.4
$ws.Range("L122").Value = 69465.59999999999
$ws.Range("N122").Value = -74365.59999999999

$ws.Range("H131").Value = 1144.0312
$ws.Range("I131").Value = 799.6667
$ws.Range("J131").Value = 1179.6552
$ws.Range("K131").Value = 2399.0001
$ws.Range("L131").Value = 3538.9656
$ws.Range("M131").Value = 2640.9999
$ws.Range("N131").Value = -13618.9656

$ws.Range("H135").Value = 925.13336
$ws.Range("I135").Value = 587
$ws.Range("K135").Value = 5283
$ws.Range("M135").Value = -2748

$ws.Range("H136").Value = 4338.385
$ws.Range("I136").Value = 766.3333
$ws.Range("J136").Value = 5410
$ws.Range("K136").Value = 2298.9999
$ws.Range("L136").Value = 16230
$ws.Range("M136").Value = 2801.0001
$ws.Range("N136").Value = -26430

$ws.Range("H137").Value = 7584668.5
$ws.Range("J137").Value = 5485.706
$ws.Range("L137").Value = 16457.118
$ws.Range("N137").Value = -26657.118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 36415376
$ws.Range("I80").Value = 72770990
$ws.Range("J80").Value = 59758.715
$ws.Range("K80").Value = 72770990
$ws.Range("L80").Value = 59758.715
$ws.Range("M80").Value = -72769992
$ws.Range("N80").Value = -61754.715

$ws.Range("H83").Value = 36415376
$ws.Range("I83").Value = 72770990
$ws.Range("J83").Value = 59758.715
$ws.Range("K83").Value = 363854950
$ws.Range("L83").Value = 298793.575
$ws.Range("M83").Value = -363849958
$ws.Range("N83").Value = -308777.575

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 144800.28
$ws.Range("I40").Value = 144800.28
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 144800.28
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -144664.28
$ws.Range("N40").ClearContents()

$ws.Range("H93").Value = 23601.6
$ws.Range("I93").Value = 51500
$ws.Range("J93").Value = 5002.6665
$ws.Range("K93").Value = 51500
$ws.Range("L93").Value = 5002.6665
$ws.Range("M93").Value = -50252
$ws.Range("N93").Value = -7498.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5053104
$ws.Range("I132").Value = 3528.5
$ws.Range("J132").Value = 8773844
$ws.Range("K132").Value = 10585.5
$ws.Range("L132").Value = 26321532
$ws.Range("M132").Value = -8055.5
$ws.Range("N132").Value = -26326592
